# Add columns I (I0) and J (IF) to the worksheet, rows 1-63
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold font, border, centered alignment) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-63 for columns I (I0) and J (IF)
$data = @(
    @(8,8),
    @(7,7),
    @(8,8),
    @(5,7),
    @(8,8),
    @(8,8),
    @(8,9),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,8),
    @(8,8),
    @(6,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(11,11),
    @(5,5),
    @(8,8),
    @(6,6),
    @(4,5),
    @(6,7),
    @(3,3),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(9,9),
    @(6,6),
    @(9,9),
    @(2,3),
    @(11,11),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(3,4),
    @(7,7),
    @(8,9),
    @(7,7),
    @(2,3),
    @(8,8),
    @(7,7),
    @(9,9),
    @(9,9),
    @(6,6),
    @(6,6),
    @(5,6),
    @(6,7),
    @(7,8),
    @(5,5),
    @(7,7),
    @(5,6),
    @(5,7),
    @(5,6),
    @(1,2)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
